$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 118
$ws1.Range("F3").Value = 508
$ws1.Range("F4").Value = 1496
$ws1.Range("F6").Value = 139
$ws1.Range("F8").Value = 5211
$ws1.Range("F9").Value = 135
$ws1.Range("F10").Value = 724
$ws1.Range("F11").Value = 1040
$ws1.Range("F12").Value = 60
$ws1.Range("F13").Value = 316
$ws1.Range("F14").Value = 46
$ws1.Range("F15").Value = 6325
$ws1.Range("F20").Value = 15166
$ws1.Range("F21").Value = 1503
$ws1.Range("F22").Value = 270
$ws1.Range("F23").Value = 131
$ws1.Range("F24").Value = 98
$ws1.Range("F25").Value = 10989
$ws1.Range("F26").Value = 734
$ws1.Range("F27").Value = 4286
$ws1.Range("F28").Value = 223
$ws1.Range("F31").Value = 124

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 118
$ws4.Range("F3").Value = 508
$ws4.Range("F4").Value = 1496
$ws4.Range("F6").Value = 139
$ws4.Range("F9").Value = 5211
$ws4.Range("F10").Value = 135
$ws4.Range("F11").Value = 724
$ws4.Range("F13").Value = 1040
$ws4.Range("F14").Value = 60
$ws4.Range("F15").Value = 316
$ws4.Range("F16").Value = 46
$ws4.Range("F18").Value = 6325
$ws4.Range("F23").Value = 15167
$ws4.Range("F24").Value = 1503
$ws4.Range("F25").Value = 270
$ws4.Range("F26").Value = 131
$ws4.Range("F27").Value = 98
$ws4.Range("F28").Value = 10989
$ws4.Range("F29").Value = 734
$ws4.Range("F30").Value = 4286
$ws4.Range("F31").Value = 223
$ws4.Range("F34").Value = 124
